$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row (row 11): correct the per-question marks
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# "Total" row (row 12): recompute totals from the corrected marking
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "46 / 112"
